# "Subida con archivo modificado a rama 2"
# Update the text of the single populated cell (A1 on Hoja1) to the new
# wording, preserving the trailing space from the source edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "My first project on excel 2 "

# The author's re-save also dropped the lingering A2 selection that had
# been left over from editing; move the active selection back onto the
# sheet's only cell (A1) to reflect that.
$ws.Range("A1").Select()
